$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 23:22"

# Update Cataluña row (row 5) figures
$ws.Range("B5").Value = 42610
$ws.Range("C5").Value = 20881
$ws.Range("D5").Value = 17153
$ws.Range("E5").Value = 4576
